# Update Wnt3-Ryk NATMI output with newly recomputed TPM-based values.
# The "ECs" sending-cluster rows are dropped from the table (no longer
# produced by the new run), and the remaining "FAPs" sending-cluster rows
# are recomputed against the smaller (FAPs-only) background.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the last three rows (old rows 5-7, the second "FAPs" sender
#    block) - their content becomes the new rows 2-4 below, so remove the
#    now-duplicate trailing rows first.
$ws.Range("A5:T7").EntireRow.Delete()

# 2) Row 2 : FAPs -> Wnt3 -> Ryk -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt3"
$ws.Range("C2").Value = "Ryk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1741663333333333
$ws.Range("H2").Value = 0.5224989999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.423863
$ws.Range("N2").Value = 22.271589
$ws.Range("O2").Value = 0.1690720838224332
$ws.Range("P2").Value = 0.1690720838224332
$ws.Range("Q2").Value = 1.292986997879
$ws.Range("R2").Value = 11.636882980911
$ws.Range("S2").Value = 0.1690720838224332
$ws.Range("T2").Value = 0.1690720838224332

# 3) Row 3 : FAPs -> Wnt3 -> Ryk -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt3"
$ws.Range("C3").Value = "Ryk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1741663333333333
$ws.Range("H3").Value = 0.5224989999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 21.552384
$ws.Range("N3").Value = 64.657152
$ws.Range("O3").Value = 0.4908369772207905
$ws.Range("P3").Value = 0.4908369772207905
$ws.Range("Q3").Value = 3.753699695872
$ws.Range("R3").Value = 33.783297262848
$ws.Range("S3").Value = 0.4908369772207905
$ws.Range("T3").Value = 0.4908369772207905

# 4) Row 4 : FAPs -> Wnt3 -> Ryk -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt3"
$ws.Range("C4").Value = "Ryk"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1741663333333333
$ws.Range("H4").Value = 0.5224989999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.93320766666667
$ws.Range("N4").Value = 44.799623
$ws.Range("O4").Value = 0.3400909389567762
$ws.Range("P4").Value = 0.3400909389567762
$ws.Range("Q4").Value = 2.600862024208555
$ws.Range("R4").Value = 23.407758217877
$ws.Range("S4").Value = 0.3400909389567762
$ws.Range("T4").Value = 0.3400909389567762
